# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 504
    $ws.Range("F3").Value = 3395
    $ws.Range("F4").Value = 92
    $ws.Range("F5").Value = 671
}
